$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "basketball leg sleeve youth padded",
    "knee pad gym",
    "work need pads",
    "softball material",
    "compression spandex men",
    "football leggings men",
    "knees bees",
    "snowboarding pants youth",
    "knee length pants",
    "medias de basketball",
    "padded sliding shorts youth",
    "men gym tight pants",
    "compression tight leggings",
    "hockey leggings girls",
    "boys black tight pants",
    "yoga pad",
    "leg compression pants men",
    "wrestling shorts men",
    "yoga pants men big and tall",
    "wrestling knee",
    "youth sliding shorts",
    "youth hockey pants",
    "sliding shorts baseball",
    "knee pads replacement",
    "thread protector 1/2 x 28",
    "running tights youth",
    "mesh leggings men",
    "boys hockey pants",
    "baseball pants youth boys black",
    "basketball tight shorts for boys",
    "sports compression leggings",
    "knee pads for yoga",
    "padded work pants mens",
    "baseball pants for boys",
    "youth spandex",
    "women knee pads for work",
    "big and tall mens compression pants",
    "padded compression shorts youth",
    "boys tights youth",
    "men sheer pants",
    "football girdle with pads for men",
    "knee pads for work men",
    "yoga position chart",
    "youth football leggings boys",
    "knee compression sleeve pad",
    "black baseball pants mens",
    "mens leggings shorts",
    "youth baseball pants knee high",
    "elastic waist baseball pants",
    "adult tights",
    "knee pads for work black",
    "arthritis test",
    "cycling sweat guard",
    "volleyball knee pads extra large",
    "size 5 basketball",
    "athletic knee compression",
    "knee pads for man",
    "youth basketball knee sleeve",
    "knee pads for mountain biking",
    "best knee pads for work",
    "hip pads for men",
    "girl compression pants",
    "compression shorts basketball",
    "basketball knee sleeve boys",
    "mens pants big and tall",
    "below knee shorts men",
    "youth wrestling shorts",
    "capri spandex",
    "yoga pads for hands",
    "football pants adult with pads",
    "boys youth compression pants",
    "exercise kneeling pad",
    "knee pad volleyball",
    "knee pads working",
    "baseball softball pants",
    "boys knee pads volleyball",
    "sliding workout pads",
    "knee pads for youth",
    "black legging for men",
    "cheap leggings for men",
    "little boys compression leggings",
    "basketball compression gear",
    "men sport pants",
    "cold knee pad",
    "black compression shorts for men",
    "calf sleeves for men basketball",
    "mens work knee pads",
    "the bees knees",
    "athletic leggings for men",
    "compression pants for boys",
    "patella knee pads",
    "lacrosse tights",
    "boys leggings youth",
    "boy sport tights",
    "volleyball pants",
    "knee pads for wrestling",
    "football padded shorts for men",
    "mens basketball outdoor",
    "compression pants size",
    "mens knee pads construction"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
